$wb = $excel.ActiveWorkbook

# "Item matches as list": append the newly-matched time-slot rows
# (Hora, Cambio, Nuevos, Actualizados) to each "22 - 9 *" sheet, right
# after the existing data rows.
function Add-Rows($sheetName, $rows) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count
    $r = $lastRow + 1
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row.Hora
        $ws.Cells.Item($r, 2).Value = $row.Cambio
        $ws.Cells.Item($r, 3).Value = $row.Nuevos
        $ws.Cells.Item($r, 4).Value = $row.Actualizados
        $r = $r + 1
    }
}

# Rows common to all four sheets: the first three new time slots never
# matched/changed anything.
$commonRows = @(
    @{ Hora = "11:28"; Cambio = $false; Nuevos = 0; Actualizados = 0 },
    @{ Hora = "15:19"; Cambio = $false; Nuevos = 0; Actualizados = 0 },
    @{ Hora = "16:45"; Cambio = $false; Nuevos = 0; Actualizados = 0 }
)

# "22 - 9 Bershka": 21:28 matched as changed, 40 items updated
Add-Rows "22 - 9 Bershka" ($commonRows + @(
    @{ Hora = "21:28"; Cambio = $true; Nuevos = 0; Actualizados = 40 }
))

# "22 - 9 Mango": 21:28 did not match/change anything
Add-Rows "22 - 9 Mango" ($commonRows + @(
    @{ Hora = "21:28"; Cambio = $false; Nuevos = 0; Actualizados = 0 }
))

# "22 - 9 Zara": 21:28 did not match/change anything
Add-Rows "22 - 9 Zara" ($commonRows + @(
    @{ Hora = "21:28"; Cambio = $false; Nuevos = 0; Actualizados = 0 }
))

# "22 - 9 Stradivarius": 21:28 matched as changed, 1 item updated
Add-Rows "22 - 9 Stradivarius" ($commonRows + @(
    @{ Hora = "21:28"; Cambio = $true; Nuevos = 0; Actualizados = 1 }
))
